$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the bottom border across into the new column Q ---
# Q3 should pick up the same style as the rest of row 3 (s=9), so copy
# formatting from the neighboring P3 cell instead of re-building borders by hand.
$ws.Range("P3").Copy($ws.Range("Q3"))

# --- Row 4: new year column header (2020) ---
# Base formatting matches P4 (year header style) but with vertical="top" alignment,
# which the engine will fold into a freshly appended cellXf.
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020
$ws.Range("Q4").VerticalAlignment = -4160

# --- Row 5: new data value (1.1) for 2020, reusing the existing row style ---
$ws.Range("Q5").Value = 1.1000000000000001

# --- Row 6: new data value (7) for 2020, with a numeric "0.0" format ---
# Base formatting on D3 (fontId=3 / borderId=1 / vertical=center, no wrap),
# then apply the "0.0" number format, which appends another new cellXf.
$ws.Range("D3").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 7
$ws.Range("Q6").NumberFormat = "0.0"

# --- Selection moves to J22 ---
$ws.Range("J22").Select()
